$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 2 (SeniorCitizen), shifting rows 3-5 up to 2-4
$ws.Rows.Item(2).Delete()

# Renumber the index column (A) to be 0-based sequential again
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
